# The presentation ships with two themes:
#   ppt/theme/theme1.xml -> bound to the (only) slide master, name="Integral"
#   ppt/theme/theme2.xml -> bound to the notes master,      name="Office Theme"
#
# The commit swaps their contents: the slide master's theme becomes the
# "Office Theme" palette (what used to live in theme2.xml) and the notes
# master's theme becomes the "Integral" palette (what used to live in
# theme1.xml). Font scheme / format scheme are identical between the two
# themes, so only the 12 colour-scheme slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) actually change.
#
# PowerPoint's COM RGB values are packed as 0x00BBGGRR (little endian
# R,G,B), so build them from the hex triples with a small helper instead
# of hand-computing each integer.
function ComRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

function HexToComRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ComRGB $r $g $b
}

# Target palette for the slide master's theme (theme1.xml): "Office Theme".
$officeTheme = [ordered]@{
    1  = "000000"  # dk1
    2  = "FFFFFF"  # lt1
    3  = "44546A"  # dk2
    4  = "E7E6E6"  # lt2
    5  = "5B9BD5"  # accent1
    6  = "ED7D31"  # accent2
    7  = "A5A5A5"  # accent3
    8  = "FFC000"  # accent4
    9  = "4472C4"  # accent5
    10 = "70AD47"  # accent6
    11 = "0563C1"  # hlink
    12 = "954F72"  # folHlink
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$tcs = $master.Theme.ThemeColorScheme

foreach ($idx in $officeTheme.Keys) {
    $tcs.Item($idx).RGB = HexToComRGB $officeTheme[$idx]
}
